# Generate Report for Handback
#
# The latest handback for bf0e4b89-b01f-44b4-844f-7e827d1b9628 came back
# pointing at a stale commit, so the report grows a "Latest Target File"
# (column I) hyperlink + an error explanation in column P for that row, on
# both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$message = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efc77b4f2d07eb1c914582117ef4b1203dfe8455/e2e/bf0e4b89-b01f-44b4-844f-7e827d1b9628.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/befe8356f8d19a35f9fc27cbfa413586829a121f/e2e/bf0e4b89-b01f-44b4-844f-7e827d1b9628.md."

# ---- zh-cn sheet, row 7 (bf0e4b89-b01f-44b4-844f-7e827d1b9628) ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/befe8356f8d19a35f9fc27cbfa413586829a121f/e2e/bf0e4b89-b01f-44b4-844f-7e827d1b9628.md",
    $null,
    $null,
    "bf0e4b89-b01f-44b4-844f-7e827d1b9628.md"
)
$wsZh.Range("J7").Value = "bf0e4b89-b01f-44b4-844f-7e827d1b9628.ae9e61f1d013a3d23faa7c21614029b3f6149d37.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-06 05:09:44"
$wsZh.Range("P7").Value = $message

# ---- de-de sheet, row 7 (bf0e4b89-b01f-44b4-844f-7e827d1b9628) ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/befe8356f8d19a35f9fc27cbfa413586829a121f/e2e/bf0e4b89-b01f-44b4-844f-7e827d1b9628.md",
    $null,
    $null,
    "bf0e4b89-b01f-44b4-844f-7e827d1b9628.md"
)
$wsDe.Range("J7").Value = "bf0e4b89-b01f-44b4-844f-7e827d1b9628.ae9e61f1d013a3d23faa7c21614029b3f6149d37.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-06 05:09:52"
$wsDe.Range("P7").Value = $message

Write-Output "Report updated for bf0e4b89-b01f-44b4-844f-7e827d1b9628 handback (zh-cn, de-de)."
